$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.59209660795975
$ws.Range("C2").Value = 10.36407734381725
$ws.Range("E2").Value = 13.68713785086328
$ws.Range("F2").Value = 54.99448763524786
$ws.Range("G2").Value = 3.697181159449026
$ws.Range("J2").Value = 10.80274813932115
$ws.Range("M2").Value = 17.01738085309923
$ws.Range("B3").Value = 17.19496184776007
$ws.Range("C3").Value = 9.982430739905219
$ws.Range("E3").Value = 13.59135228416763
$ws.Range("F3").Value = 53.9563368386496
$ws.Range("G3").Value = 3.702172325133112
$ws.Range("J3").Value = 10.69557174660276
$ws.Range("M3").Value = 17.16592048138934
$ws.Range("B4").Value = 16.95447055075539
$ws.Range("C4").Value = 9.746123393369515
$ws.Range("E4").Value = 13.53554460703709
$ws.Range("F4").Value = 53.31651361752625
$ws.Range("G4").Value = 3.705388615612043
$ws.Range("J4").Value = 10.62938246835669
$ws.Range("M4").Value = 17.26428446602803
$ws.Range("B5").Value = 16.85748143073207
$ws.Range("C5").Value = 9.649526498575762
$ws.Range("E5").Value = 13.5135670475943
$ws.Range("F5").Value = 53.05544196961425
$ws.Range("G5").Value = 3.706737605244185
$ws.Range("J5").Value = 10.60232318111156
$ws.Range("M5").Value = 17.30616457580034
$ws.Range("B6").Value = 16.84144264560834
$ws.Range("C6").Value = 9.633474600509203
$ws.Range("E6").Value = 13.50996409482794
$ws.Range("F6").Value = 53.01207854955774
$ws.Range("G6").Value = 3.706963923966562
$ws.Range("J6").Value = 10.59782505100208
$ws.Range("M6").Value = 17.31322708296837
$ws.Range("B7").Value = 16.95315819380339
$ws.Range("C7").Value = 9.744821596261161
$ws.Range("E7").Value = 13.53524510386724
$ws.Range("F7").Value = 53.3129937335879
$ws.Range("G7").Value = 3.705406653139903
$ws.Range("J7").Value = 10.62901787558664
$ws.Range("M7").Value = 17.26484200961728
$ws.Range("B8").Value = 17.45457075420742
$ws.Range("C8").Value = 10.23300881231216
$ws.Range("E8").Value = 13.65349307375761
$ws.Range("F8").Value = 54.63718705122275
$ws.Range("G8").Value = 3.698870738104855
$ws.Range("J8").Value = 10.76587890841827
$ws.Range("M8").Value = 17.06710849821127
$ws.Range("B9").Value = 18.45706258157559
$ws.Range("C9").Value = 11.16682406244751
$ws.Range("E9").Value = 13.9087626524597
$ws.Range("F9").Value = 57.20360995035947
$ws.Range("G9").Value = 3.687249211059044
$ws.Range("J9").Value = 11.03080401934798
$ws.Range("M9").Value = 16.73637625318306
$ws.Range("B10").Value = 19.19576822495313
$ws.Range("C10").Value = 11.82920862800333
$ws.Range("E10").Value = 14.1098617466245
$ws.Range("F10").Value = 59.05603789624417
$ws.Range("G10").Value = 3.679428037550921
$ws.Range("J10").Value = 11.22276786548663
$ws.Range("M10").Value = 16.52845549362621
$ws.Range("B11").Value = 19.53040735546759
$ws.Range("C11").Value = 12.12375224632919
$ws.Range("E11").Value = 14.20411341232749
$ws.Range("F11").Value = 59.88864208513613
$ws.Range("G11").Value = 3.676023214305274
$ws.Range("J11").Value = 11.30939428229841
$ws.Range("M11").Value = 16.44155321070623
$ws.Range("B12").Value = 19.65677890367034
$ws.Range("C12").Value = 12.23419735323519
$ws.Range("E12").Value = 14.24018582576866
$ws.Range("F12").Value = 60.20224918270964
$ws.Range("G12").Value = 3.674755712384335
$ws.Range("J12").Value = 11.34208651100623
$ws.Range("M12").Value = 16.40975668481488
$ws.Range("B13").Value = 19.62958030018935
$ws.Range("C13").Value = 12.21046131342043
$ws.Range("E13").Value = 14.23240030698027
$ws.Range("F13").Value = 60.1347866758048
$ws.Range("G13").Value = 3.675027723546517
$ws.Range("J13").Value = 11.33505077670496
$ws.Range("M13").Value = 16.41655509000752
$ws.Range("B14").Value = 19.54081182493101
$ws.Range("C14").Value = 12.13286114135755
$ws.Range("E14").Value = 14.20707356966821
$ws.Range("F14").Value = 59.91447773295954
$ws.Range("G14").Value = 3.675918499673255
$ws.Range("J14").Value = 11.3120861841836
$ws.Range("M14").Value = 16.43891497543713
$ws.Range("B15").Value = 19.48638890122211
$ws.Range("C15").Value = 12.08518332644951
$ws.Range("E15").Value = 14.19160935269061
$ws.Range("F15").Value = 59.77930637838305
$ws.Range("G15").Value = 3.676466963420905
$ws.Range("J15").Value = 11.29800489212791
$ws.Range("M15").Value = 16.45275600657693
$ws.Range("B16").Value = 19.17385774573984
$ws.Range("C16").Value = 11.80981273532644
$ws.Range("E16").Value = 14.10375649469969
$ws.Range("F16").Value = 59.00140230400817
$ws.Range("G16").Value = 3.679653615478549
$ws.Range("J16").Value = 11.21709163560813
$ws.Range("M16").Value = 16.53428991918427
$ws.Range("B17").Value = 18.98166312424408
$ws.Range("C17").Value = 11.63905741252093
$ws.Range("E17").Value = 14.05055941544073
$ws.Range("F17").Value = 58.52143315189523
$ws.Range("G17").Value = 3.681647599551373
$ws.Range("J17").Value = 11.16726725204195
$ws.Range("M17").Value = 16.58628062280308
$ws.Range("B18").Value = 18.8709941696013
$ws.Range("C18").Value = 11.54021149826466
$ws.Range("E18").Value = 14.02022373917105
$ws.Range("F18").Value = 58.24443253272241
$ws.Range("G18").Value = 3.682808904724408
$ws.Range("J18").Value = 11.13854444457143
$ws.Range("M18").Value = 16.61690693152116
$ws.Range("B19").Value = 18.83350696411084
$ws.Range("C19").Value = 11.50663924925478
$ws.Range("E19").Value = 14.00999808220995
$ws.Range("F19").Value = 58.15049159419139
$ws.Range("G19").Value = 3.683204584534593
$ws.Range("J19").Value = 11.12880853842236
$ws.Range("M19").Value = 16.62740038249998
$ws.Range("B20").Value = 19.00213646931024
$ws.Range("C20").Value = 11.65730091181448
$ws.Range("E20").Value = 14.05619535270436
$ws.Range("F20").Value = 58.57262519148097
$ws.Range("G20").Value = 3.681433845610935
$ws.Range("J20").Value = 11.17257796303308
$ws.Range("M20").Value = 16.58067128556079
$ws.Range("B21").Value = 19.56689587881089
$ws.Range("C21").Value = 12.15568469786444
$ws.Range("E21").Value = 14.21450243717444
$ws.Range("F21").Value = 59.97923525818104
$ws.Range("G21").Value = 3.675656266089883
$ws.Range("J21").Value = 11.31883454851372
$ws.Range("M21").Value = 16.43231711379763
$ws.Range("B22").Value = 19.9339027533002
$ws.Range("C22").Value = 12.47499689898489
$ws.Range("E22").Value = 14.32017954230936
$ws.Range("F22").Value = 60.8886326593881
$ws.Range("G22").Value = 3.672007446188867
$ws.Range("J22").Value = 11.41376790558821
$ws.Range("M22").Value = 16.34184104062199
$ws.Range("B23").Value = 19.73826232285494
$ws.Range("C23").Value = 12.30519560472955
$ws.Range("E23").Value = 14.26358094868561
$ws.Range("F23").Value = 60.40425019434583
$ws.Range("G23").Value = 3.673943314352403
$ws.Range("J23").Value = 11.36316355525212
$ws.Range("M23").Value = 16.3895343889765
$ws.Range("B24").Value = 18.99288099977051
$ws.Range("C24").Value = 11.649055123566
$ws.Range("E24").Value = 14.05364657161607
$ws.Range("F24").Value = 58.54948455408692
$ws.Range("G24").Value = 3.681530437227772
$ws.Range("J24").Value = 11.17017723305438
$ws.Range("M24").Value = 16.58320497446592
$ws.Range("B25").Value = 18.18481266201361
$ws.Range("C25").Value = 10.91774968888772
$ws.Range("E25").Value = 13.83726926136747
$ws.Range("F25").Value = 56.51411443271653
$ws.Range("G25").Value = 3.690266363506812
$ws.Range("J25").Value = 10.95956332642939
$ws.Range("M25").Value = 16.81971290508394
